$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.309.53"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.363.60"
$ws.Range("E3").Value = "  +7.58%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.38"
$ws.Range("E5").Value = "  +6.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "628.55"
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.40"
$ws.Range("E7").Value = "  +24.25%  "
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +11.48%  "
$ws.Range("D11").Value = "3.358.52"
$ws.Range("E11").Value = "  +7.45%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "98.421.10"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.42"
$ws.Range("E14").Value = "  +7.48%  "
$ws.Range("E15").Value = "  +3.91%  "
$ws.Range("D16").Value = "3.991.97"
$ws.Range("E16").Value = "  +7.62%  "
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "3.371.69"
$ws.Range("E18").Value = "  +7.57%  "
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.21"
$ws.Range("E20").Value = "  +5.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "488.67"
$ws.Range("E21").Value = "  -5.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").Value = "  +8.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000210"
$ws.Range("E23").Value = "  +10.16%  "
$ws.Range("E24").Value = "  +7.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.64"
$ws.Range("E25").Value = "  +3.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.04"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.93"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("E29").Value = "  +16.33%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +9.18%  "
$ws.Range("E32").Value = "  +10.34%  "
# Row 33 and 34 swap: InternetComputer(DFINITY) and Binance-PegBSC-USD swap places with updated values
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.63"
$ws.Range("E34").Value = "  +7.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.06"
$ws.Range("E35").Value = "  +6.04%  "
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.29"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  +4.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "499.97"
$ws.Range("E39").Value = "  +7.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.92"
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.460"
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.74"
$ws.Range("E42").Value = "  +4.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.27"
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.29"
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.794"
$ws.Range("E45").Value = "  +14.45%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.44"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.850"
$ws.Range("E49").Value = "  +12.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.63"
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("E51").Value = "  +3.87%  "
